# Updates a batch of computed market-price / profit figures across the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) following a
# scheduled data refresh (currentAveragePrice* / LevePrice* / LeveProfit*
# columns H:N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 30429780
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 36515136
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 36515136
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -36515948

$ws.Range("H91").Value = 30429780
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 36515136
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 36515136
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -36517944

$ws.Range("H116").Value = 1942.2667
$ws.Range("I116").Value = 1744
$ws.Range("J116").Value = 2487.5
$ws.Range("K116").Value = 1744
$ws.Range("L116").Value = 2487.5
$ws.Range("M116").Value = 1698
$ws.Range("N116").Value = -9371.5

$ws.Range("H129").Value = 1017.30554
$ws.Range("I129").Value = 1750.25
$ws.Range("K129").Value = 5250.75
$ws.Range("M129").Value = -250.75

$ws.Range("H132").Value = 7147503.5
$ws.Range("I132").Value = 9263392
$ws.Range("J132").Value = 6379.8125
$ws.Range("K132").Value = 27790176
$ws.Range("L132").Value = 19139.4375
$ws.Range("M132").Value = -27787646
$ws.Range("N132").Value = -24199.4375

$ws.Range("H137").Value = 1379.65
$ws.Range("J137").Value = 1750.25
$ws.Range("L137").Value = 5250.75
$ws.Range("N137").Value = -10350.75

$ws.Range("H138").Value = 1762.2
$ws.Range("I138").Value = 678.7646999999999
$ws.Range("J138").Value = 3179
$ws.Range("K138").Value = 2036.2941
$ws.Range("L138").Value = 9537
$ws.Range("M138").Value = 3103.7059
$ws.Range("N138").Value = -19817


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 44555.22
$ws.Range("I45").Value = 63479.5
$ws.Range("J45").Value = 1299.7142
$ws.Range("K45").Value = 63479.5
$ws.Range("L45").Value = 1299.7142
$ws.Range("M45").Value = -63102.5
$ws.Range("N45").Value = -2053.7142

$ws.Range("H61").Value = 1279.1052
$ws.Range("I61").Value = 1125.28
$ws.Range("J61").Value = 1574.9231
$ws.Range("K61").Value = 1125.28
$ws.Range("L61").Value = 1574.9231
$ws.Range("M61").Value = -913.28
$ws.Range("N61").Value = -1998.9231

$ws.Range("H74").Value = 837.5
$ws.Range("I74").Value = 821.79486
$ws.Range("K74").Value = 821.79486
$ws.Range("M74").Value = 52.20514000000003

$ws.Range("H77").Value = 837.5
$ws.Range("I77").Value = 821.79486
$ws.Range("K77").Value = 4108.9743
$ws.Range("M77").Value = 259.0257000000001

$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws.Range("H136").Value = 1279.1052
$ws.Range("I136").Value = 1125.28
$ws.Range("J136").Value = 1574.9231
$ws.Range("K136").Value = 3375.84
$ws.Range("L136").Value = 4724.7693
$ws.Range("M136").Value = -825.8400000000001
$ws.Range("N136").Value = -9824.7693


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2317833
$ws.Range("I134").Value = 778.5454999999999
$ws.Range("K134").Value = 2335.6365
$ws.Range("M134").Value = 199.3635000000004


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1153.6769
$ws.Range("I31").Value = 846.85
$ws.Range("K31").Value = 846.85
$ws.Range("M31").Value = -551.85

$ws.Range("H34").Value = 1153.6769
$ws.Range("I34").Value = 846.85
$ws.Range("K34").Value = 846.85
$ws.Range("M34").Value = -644.85

$ws.Range("H58").Value = 18519440
$ws.Range("J58").Value = 680.86664
$ws.Range("L58").Value = 680.86664
$ws.Range("N58").Value = -1086.86664

$ws.Range("H99").Value = 30306152
$ws.Range("I99").Value = 50002710
$ws.Range("J99").Value = 3753.8462
$ws.Range("K99").Value = 50002710
$ws.Range("L99").Value = 3753.8462
$ws.Range("M99").Value = -50001212
$ws.Range("N99").Value = -6749.8462

$ws.Range("H126").Value = 30306152
$ws.Range("I126").Value = 50002710
$ws.Range("J126").Value = 3753.8462
$ws.Range("K126").Value = 150008130
$ws.Range("L126").Value = 11261.5386
$ws.Range("M126").Value = -150005660
$ws.Range("N126").Value = -16201.5386

$ws.Range("H132").Value = 9260629
$ws.Range("I132").Value = 1225.1428
$ws.Range("J132").Value = 41668544
$ws.Range("K132").Value = 3675.4284
$ws.Range("L132").Value = 125005632
$ws.Range("M132").Value = -1145.4284
$ws.Range("N132").Value = -125010692

$ws.Range("H136").Value = 18519440
$ws.Range("J136").Value = 680.86664
$ws.Range("L136").Value = 2042.59992
$ws.Range("N136").Value = -7142.59992


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 43860700
$ws.Range("I137").Value = 35715268
$ws.Range("J137").Value = 66667908
$ws.Range("K137").Value = 107145804
$ws.Range("L137").Value = 200003724
$ws.Range("M137").Value = -107140704
$ws.Range("N137").Value = -200013924

$ws.Range("H141").Value = 43480816
$ws.Range("I141").Value = 55557540
$ws.Range("J141").Value = 4598.8
$ws.Range("K141").Value = 166672620
$ws.Range("L141").Value = 13796.4
$ws.Range("M141").Value = -166667440
$ws.Range("N141").Value = -24156.4


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1211.6111
$ws.Range("I113").Value = 1145.1111
$ws.Range("J113").Value = 1278.1111
$ws.Range("K113").Value = 1145.1111
$ws.Range("L113").Value = 1278.1111
$ws.Range("M113").Value = 1024.8889
$ws.Range("N113").Value = -5618.1111

$ws.Range("H132").Value = 4703.1836
$ws.Range("I132").Value = 2850.1333
$ws.Range("K132").Value = 8550.3999
$ws.Range("M132").Value = -6020.3999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9617650
$ws.Range("I40").Value = 1701.6471
$ws.Range("J40").Value = 27781106
$ws.Range("K40").Value = 1701.6471
$ws.Range("L40").Value = 27781106
$ws.Range("M40").Value = -1565.6471
$ws.Range("N40").Value = -27781378

$ws.Range("H46").Value = 2852.8333
$ws.Range("I46").Value = 711.6667
$ws.Range("J46").Value = 4994
$ws.Range("K46").Value = 711.6667
$ws.Range("L46").Value = 4994
$ws.Range("M46").Value = -523.6667
$ws.Range("N46").Value = -5370

$ws.Range("H61").Value = 2025.0625
$ws.Range("I61").Value = 1984.6923
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 1984.6923
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -1782.6923
$ws.Range("N61").Value = -2604

$ws.Range("H113").Value = 2025.0625
$ws.Range("I113").Value = 1984.6923
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1984.6923
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 185.3077000000001
$ws.Range("N113").Value = -6540

$ws.Range("H136").Value = 30889910
$ws.Range("I136").Value = 4083475.5
$ws.Range("J136").Value = 500002500
$ws.Range("K136").Value = 12250426.5
$ws.Range("L136").Value = 1500007500
$ws.Range("M136").Value = -12247876.5
$ws.Range("N136").Value = -1500012600


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 40000330
$ws.Range("I113").Value = 55555868
$ws.Range("J113").Value = 370.42856
$ws.Range("K113").Value = 166667604
$ws.Range("L113").Value = 1111.28568
$ws.Range("M113").Value = -166665434
$ws.Range("N113").Value = -5451.28568

$ws.Range("H122").Value = 15148.405
$ws.Range("I122").Value = 22013.125
$ws.Range("J122").Value = 2475.077
$ws.Range("K122").Value = 66039.375
$ws.Range("L122").Value = 7425.231000000001
$ws.Range("M122").Value = -63589.375
$ws.Range("N122").Value = -12325.231

$ws.Range("H126").Value = 1219.8572
$ws.Range("J126").Value = 2200
$ws.Range("L126").Value = 6600
$ws.Range("N126").Value = -11540

$ws.Range("H132").Value = 20954.701
$ws.Range("I132").Value = 22175.166
$ws.Range("J132").Value = 14445.556
$ws.Range("K132").Value = 66525.49800000001
$ws.Range("L132").Value = 43336.66800000001
$ws.Range("M132").Value = -63995.49800000001
$ws.Range("N132").Value = -48396.66800000001
